$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove excess "NO. OF OVERTIME HOURS" (column I) entries for rows that
#    should no longer carry a value, and populate the "NO. OF HOURS
#    UNDERTIME" (column F) for rows 7 and 10.
# ---------------------------------------------------------------------------
$ws.Range("F7").Value = 5.5
$ws.Range("I7").ClearContents()

$ws.Range("I8").ClearContents()

$ws.Range("I9").ClearContents()

$ws.Range("F10").Value = 2.75
$ws.Range("I10").ClearContents()

$ws.Range("I14").ClearContents()
$ws.Range("I15").ClearContents()
$ws.Range("I16").ClearContents()
$ws.Range("I17").ClearContents()
$ws.Range("I18").ClearContents()

# ---------------------------------------------------------------------------
# 2. Add a "Legends:" section under the summary table (rows 24-30).
# ---------------------------------------------------------------------------

# "Legends:" header, styled like the other big bold/underlined headers.
$legendHeader = $ws.Range("E24:P24")
$legendHeader.Font.Name = "Arial"
$legendHeader.Font.Size = 15
$legendHeader.Font.Bold = $true
$legendHeader.Font.Underline = $true
$legendHeader.Merge()
$ws.Range("E24").Value = "Legends:"

# Legend entry 1 - blue swatch / remarks note.
$swatch1 = $ws.Range("E25:E26")
$swatch1.Interior.Color = 13411113
$swatch1.Merge()

$note1 = $ws.Range("F25:P26")
$note1.Font.Name = "Arial"
$note1.Font.Size = 11
$note1.Font.Bold = $true
$note1.Font.Underline = $true
$note1.Merge()
$ws.Range("F25").Value = "Employee has request(s)/remark(s) for that day." + [char]10 + "*May incur late and/or undertime depending on his or her time-in and time-out."

# Legend entry 2 - orange swatch / half-day note.
$swatch2 = $ws.Range("E27:E28")
$swatch2.Interior.Color = 6737151
$swatch2.Merge()

$note2 = $ws.Range("F27:P28")
$note2.Font.Name = "Arial"
$note2.Font.Size = 11
$note2.Font.Bold = $true
$note2.Font.Underline = $true
$note2.Merge()
$ws.Range("F27").Value = "Employee is considered half-day because of his time-in or time-out."

# Legend entry 3 - red swatch / absent note.
$swatch3 = $ws.Range("E29:E30")
$swatch3.Interior.Color = 6184671
$swatch3.Merge()

$note3 = $ws.Range("F29:P30")
$note3.Font.Name = "Arial"
$note3.Font.Size = 11
$note3.Font.Bold = $true
$note3.Font.Underline = $true
$note3.Merge()
$ws.Range("F29").Value = "Employee has no time-in and therefore, considered as absent."
